# Generate Report for Handoff
# Refresh the localization-status report: the three tracked files
# (707a38ce..., ffff72dcf5a8..., ffffff4b4f9267...) get reshuffled across
# rows 2-4 on every sheet, row 4's status flips from "Handed back: in sync
# with en-US" to "Ready for handoff" (with matching new handoff/handback
# file names, dates and an out-of-date error message), and the Content
# Duplicate flags for rows 3/4 swap.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md"
$ws1.Range("A3").Value = "ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md"
$ws1.Range("A4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md"

$ws1.Range("B2").Value = "e2e\ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md"
$ws1.Range("B3").Value = "e2e\ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md"
$ws1.Range("B4").Value = "e2e\707a38ce-5d43-4dd2-818b-5c55153ba10b.md"

foreach ($hl in $ws1.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') { $hl.TextToDisplay = "e2e\ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md" }
    if ($addr -eq '$B$3') { $hl.TextToDisplay = "e2e\ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md" }
    if ($addr -eq '$B$4') { $hl.TextToDisplay = "e2e\707a38ce-5d43-4dd2-818b-5c55153ba10b.md" }
}

$ws1.Range("G2").Value = "2016-08-17 15:05:27"

$ws1.Range("E4").Value = "Ready for handoff"
$ws1.Range("F4").Value = "Ready for handoff"
$ws1.Range("G4").Value = "2016-08-17 15:07:50"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md"
$ws2.Range("A3").Value = "ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md"
$ws2.Range("A4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md"

foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md" }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md" }
    if ($addr -eq '$A$4') { $hl.TextToDisplay = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md" }
}

$ws2.Range("C4").Value = "Ready for handoff"

# F3/F4 ("Content Duplicate") swap True/False. Assigning the literal text
# "True"/"False" via .Value gets auto-coerced to an Excel boolean cell
# (t="b"), but the sheet stores these as plain shared-string text (t="s")
# elsewhere (e.g. F2), so swap the two cells with Copy instead, which
# preserves the original string cell type.
$ws2.Range("F4").Copy($ws2.Range("Z1"))
$ws2.Range("F3").Copy($ws2.Range("F4"))
$ws2.Range("Z1").Copy($ws2.Range("F3"))
$ws2.Range("Z1").ClearContents()

$ws2.Range("G2").Value = "893f6d0d-6efc-4983-a846-aa6fe86977c7.2c11a8e109289ee1c6b619a9a637f49868cff920.zh-cn.xlf"
$ws2.Range("G4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.f9c6f63c30a28391b56341230c2005b3bc868e93.zh-cn.xlf"

$ws2.Range("H2").Value = "2016-08-17 15:05:22"
$ws2.Range("H4").Value = "2016-08-17 15:07:44"

$ws2.Range("I2").Value = "893f6d0d-6efc-4983-a846-aa6fe86977c7.md"
$ws2.Range("I4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md"

foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') { $hl.TextToDisplay = "893f6d0d-6efc-4983-a846-aa6fe86977c7.md" }
    if ($addr -eq '$I$4') { $hl.TextToDisplay = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md" }
}

$ws2.Range("J2").Value = "893f6d0d-6efc-4983-a846-aa6fe86977c7.2c11a8e109289ee1c6b619a9a637f49868cff920.zh-cn.xlf"
$ws2.Range("J4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.f9c6f63c30a28391b56341230c2005b3bc868e93.zh-cn.xlf"

$ws2.Range("K2").Value = "2016-08-17 15:05:49"
$ws2.Range("K4").Value = "2016-08-17 15:07:19"

$ws2.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a99b7e83c402b953edaab33a38f47ab8b7b830a/e2e/707a38ce-5d43-4dd2-818b-5c55153ba10b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce24c6bfbeb5596a827a5c98ff0d5c0cde9c01d9/e2e/707a38ce-5d43-4dd2-818b-5c55153ba10b.md."

# Column P widened to fit the new error message.
$ws2.Columns.Item(16).ColumnWidth = $ws2.Columns.Item(1).ColumnWidth

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md"
$ws3.Range("A3").Value = "ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md"
$ws3.Range("A4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md"

foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff72dcf5a8-302c-44c1-979c-79c7da3ae30e.md" }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffff4b4f9267-3eaf-4b84-8c9b-9c731aaf673d.md" }
    if ($addr -eq '$A$4') { $hl.TextToDisplay = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md" }
}

$ws3.Range("C4").Value = "Ready for handoff"

# Same True/False swap as on the "zh-cn" sheet; see comment there.
$ws3.Range("F4").Copy($ws3.Range("Z1"))
$ws3.Range("F3").Copy($ws3.Range("F4"))
$ws3.Range("Z1").Copy($ws3.Range("F3"))
$ws3.Range("Z1").ClearContents()

$ws3.Range("G2").Value = "893f6d0d-6efc-4983-a846-aa6fe86977c7.2c11a8e109289ee1c6b619a9a637f49868cff920.de-de.xlf"
$ws3.Range("G4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.f9c6f63c30a28391b56341230c2005b3bc868e93.de-de.xlf"

$ws3.Range("H2").Value = "2016-08-17 15:05:27"
$ws3.Range("H4").Value = "2016-08-17 15:07:50"

$ws3.Range("I2").Value = "893f6d0d-6efc-4983-a846-aa6fe86977c7.md"
$ws3.Range("I4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md"

foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') { $hl.TextToDisplay = "893f6d0d-6efc-4983-a846-aa6fe86977c7.md" }
    if ($addr -eq '$I$4') { $hl.TextToDisplay = "707a38ce-5d43-4dd2-818b-5c55153ba10b.md" }
}

$ws3.Range("J2").Value = "893f6d0d-6efc-4983-a846-aa6fe86977c7.2c11a8e109289ee1c6b619a9a637f49868cff920.de-de.xlf"
$ws3.Range("J4").Value = "707a38ce-5d43-4dd2-818b-5c55153ba10b.f9c6f63c30a28391b56341230c2005b3bc868e93.de-de.xlf"

$ws3.Range("K2").Value = "2016-08-17 15:05:57"
$ws3.Range("K4").Value = "2016-08-17 15:07:27"

$ws3.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a99b7e83c402b953edaab33a38f47ab8b7b830a/e2e/707a38ce-5d43-4dd2-818b-5c55153ba10b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce24c6bfbeb5596a827a5c98ff0d5c0cde9c01d9/e2e/707a38ce-5d43-4dd2-818b-5c55153ba10b.md."

# Column P widened to fit the new error message.
$ws3.Columns.Item(16).ColumnWidth = $ws3.Columns.Item(1).ColumnWidth
